# Update "想去人数" (number of people interested) figures that changed
# between site generations.
#
# Sheet "展览" (sheet1 / first tab) and sheet "全部类型" (sheet4 / fourth tab)
# both list the same events (in slightly different row orders because the
# latter sheet has one extra row), so the same F-column values need to be
# bumped on each of them.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3161
$ws1.Range("F3").Value = 732
$ws1.Range("F4").Value = 116
$ws1.Range("F5").Value = 6859
$ws1.Range("F6").Value = 1951
$ws1.Range("F7").Value = 15
$ws1.Range("F14").Value = 172

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3161
$ws4.Range("F4").Value = 732
$ws4.Range("F5").Value = 116
$ws4.Range("F6").Value = 6859
$ws4.Range("F7").Value = 1951
$ws4.Range("F8").Value = 15
$ws4.Range("F15").Value = 172
